# Add another survey response (row 38) to the dataset.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$row = 38

# Numeric / date cells first.
$ws.Cells.Item($row, 1).Value = 37
$ws.Cells.Item($row, 3).Value = 850
$ws.Cells.Item($row, 4).Value = 225
$ws.Cells.Item($row, 9).Value = 2

# Text cells that reuse already-existing shared strings.
$ws.Cells.Item($row, 2).Value = "Arthur"
$ws.Cells.Item($row, 5).Value = "soybean,corn silage"
$ws.Cells.Item($row, 8).Value = "yes"
$ws.Cells.Item($row, 10).Value = "winter rye"
$ws.Cells.Item($row, 11).Value = "small grains"

# Text cells that introduce brand-new shared strings -- set in the same
# order the original authoring tool wrote them so new shared-string table
# entries land at the same indices as the target file.
$ws.Cells.Item($row, 12).Value = "110"
$ws.Cells.Item($row, 15).Value = "solid"
$ws.Cells.Item($row, 18).Value = "before,after"
$ws.Cells.Item($row, 16).Value = "24,24"
$ws.Cells.Item($row, 17).Value = "solid,solid"
$ws.Cells.Item($row, 24).Value = "soil conservation,nutrient management,manure management,nutrient cycling,incentive payment,public relations"
$ws.Cells.Item($row, 26).Value = "termination economics"

# Remaining text cells that reuse already-existing shared strings.
$ws.Cells.Item($row, 13).Value = "no-till drill"
$ws.Cells.Item($row, 14).Value = "no-till"
$ws.Cells.Item($row, 19).Value = "bin run"
$ws.Cells.Item($row, 20).Value = "yes"
$ws.Cells.Item($row, 21).Value = "yes"
$ws.Cells.Item($row, 22).Value = "cover only then plant,harvest grain then plant"
$ws.Cells.Item($row, 23).Value = "plant green"
$ws.Cells.Item($row, 25).Value = "late harvest"

# Apply the existing date number format (style index 1) to the two date
# cells *before* writing their values, so Excel reuses that style instead
# of minting a brand new (orphan) numFmt/xf pair.
$ws.Cells.Item(2, 6).Copy()
$ws.Cells.Item($row, 6).PasteSpecial(-4122)  # xlPasteFormats
$ws.Cells.Item(2, 7).Copy()
$ws.Cells.Item($row, 7).PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Cells.Item($row, 6).Value = "10/15/2022"
$ws.Cells.Item($row, 7).Value = "11/8/2022"

# Update the selection to reflect where the author was working when the
# new row was added.
$ws.Activate()
$ws.Range("K19").Select()
